# Weekly fruit/vegetable price update: a new record is inserted as row 603
# (pushing the existing rows 603-648 down to 604-649) on the single data
# sheet of this "Uva" (grape) price workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 603; Excel shifts rows
# 603-648 down to 604-649 and carries their formatting (and the rows
# below keep all of their original values untouched).
$ws.Rows(603).Insert()

# Populate the newly inserted row 603 with the new weekly price record.
$ws.Cells.Item(603, 1).Value  = 9
$ws.Cells.Item(603, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(603, 3).Value  = "Metropolitana"
$ws.Cells.Item(603, 4).Value  = 44746
$ws.Cells.Item(603, 5).Value  = 13
$ws.Cells.Item(603, 6).Value  = "Fruta"
$ws.Cells.Item(603, 7).Value  = 100109
$ws.Cells.Item(603, 8).Value  = "Uva"
$ws.Cells.Item(603, 9).Value  = 100109001
$ws.Cells.Item(603, 10).Value = "Uva"
$ws.Cells.Item(603, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(603, 12).Value = "Primera"
$ws.Cells.Item(603, 13).Value = 200
$ws.Cells.Item(603, 14).Value = 6000
$ws.Cells.Item(603, 15).Value = 6000
$ws.Cells.Item(603, 16).Value = 6000
$ws.Cells.Item(603, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(603, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(603, 19).Value = 750
$ws.Cells.Item(603, 20).Value = 8
